$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("محاسبات راندمان")

# Fuel composition (row 2) - set to pure CH4 (100%), everything else 0
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# N2 formula: sum now includes L2
$ws.Range("N2").Formula = "=SUM(B2:L2)"

# Ambient temperature (B11) and flue/chimney temperature (B12)
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 250

# O2 percent (B16): burner number changed 3 -> 6
$ws.Range("B16").Value = 6

# New helper formulas
$ws.Range("Q3").Formula = "=I3*I5+K3*K5+L3*L5"
$ws.Range("Q4").Formula = "=Q3/N9"
$ws.Range("O19").Formula = "=(N17*B17*O5)/N5"
